$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New exercise rows: column A = short label, column B = github link
$labels = @(
    "Luyện tập vòng lặp 01 - Bài 01",
    "Luyện tập vòng lặp 01 - Bài 02",
    "Luyện tập vòng lặp 01 - Bài 03",
    "Luyện tập vòng lặp 01 - Bài 04",
    "Luyện tập vòng lặp 01 - Bài 05",
    "Luyện tập vòng lặp 01 - Bài 06",
    "Luyện tập vòng lặp 01 - Bài 07",
    "Luyện tập vòng lặp 01 - Bài 08"
)

$links = @(
    "https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s11_vong%20lap%202/%5BB%C3%A0i%20t%E1%BA%ADp%201%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20v%C3%B2ng%20l%E1%BA%B7p%2001.html",
    "https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s11_vong%20lap%202/%5BB%C3%A0i%20t%E1%BA%ADp%202%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20v%C3%B2ng%20l%E1%BA%B7p%2001.html",
    "https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s11_vong%20lap%202/%5BB%C3%A0i%20t%E1%BA%ADp%203%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20v%C3%B2ng%20l%E1%BA%B7p%2001.html",
    "https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s11_vong%20lap%202/%5BB%C3%A0i%20t%E1%BA%ADp%204%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20v%C3%B2ng%20l%E1%BA%B7p%2001.html",
    "https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s11_vong%20lap%202/%5BB%C3%A0i%20t%E1%BA%ADp%205%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20v%C3%B2ng%20l%E1%BA%B7p%2001.html",
    "https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s11_vong%20lap%202/%5BB%C3%A0i%20t%E1%BA%ADp%206%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20v%C3%B2ng%20l%E1%BA%B7p%2001.html",
    "https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s11_vong%20lap%202/%5BB%C3%A0i%20t%E1%BA%ADp%207%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20v%C3%B2ng%20l%E1%BA%B7p%2001.html",
    "https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s11_vong%20lap%202/%5BB%C3%A0i%20t%E1%BA%ADp%208%5D%20Luy%E1%BB%87n%20t%E1%BA%ADp%20v%C3%B2ng%20l%E1%BA%B7p%2001.html"
)

$startRow = 39
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}
for ($i = 0; $i -lt $links.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $links[$i]
}

# Update the view to reflect the new content having been scrolled to / selected
$ws.Range("A48").Select()
